$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value2 = 330
$ws.Range("F2").Value2 = 45078
$ws.Range("E5").Value2 = 29891
$ws.Range("F5").Value2 = 45078
$ws.Range("E7").Value2 = 29860
$ws.Range("F7").Value2 = 45047
$ws.Range("E9").Value2 = 29860
$ws.Range("F9").Value2 = 45047
$ws.Range("C14").Value2 = 419
$ws.Range("F14").Value2 = 45047
$ws.Range("C16").Value2 = 450
$ws.Range("F16").Value2 = 45047
$ws.Range("C17").Value2 = 365
$ws.Range("F17").Value2 = 45047
$ws.Range("C21").Value2 = 305
$ws.Range("F21").Value2 = 45047
$ws.Range("C23").Value2 = 255
$ws.Range("F23").Value2 = 45017
$ws.Range("C24").Value2 = 211
$ws.Range("F24").Value2 = 45078
$ws.Range("E30").Value2 = 29860
$ws.Range("F30").Value2 = 45047
$ws.Range("C31").Value2 = 389
$ws.Range("F31").Value2 = 45047
$ws.Range("C33").Value2 = 450
$ws.Range("F33").Value2 = 45078
$ws.Range("C34").Value2 = 197
$ws.Range("F34").Value2 = 45047
$ws.Range("C36").Value2 = 498
$ws.Range("F36").Value2 = 45078
$ws.Range("C41").Value2 = 389
$ws.Range("F41").Value2 = 45047
$ws.Range("C42").Value2 = 234
$ws.Range("F42").Value2 = 45047
$ws.Range("C43").Value2 = 468
$ws.Range("F43").Value2 = 45047
$ws.Range("C44").Value2 = 401
$ws.Range("F44").Value2 = 45047
$ws.Range("C48").Value2 = 353
$ws.Range("F48").Value2 = 45047
$ws.Range("C52").Value2 = 342
$ws.Range("F52").Value2 = 45078
